# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "In Translation" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File / Latest Handback DateTime columns
#    are populated (with hyperlinks) for both the zh-cn and de-de language sheets
#  - Column widths are widened to fit the new, longer content

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$docUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cd994b6632ff6c5449546a3e51fe36981779d8b/e2e/e25f3909-71f0-446d-a71c-4e2ddb72b79a.md"
$docUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cd994b6632ff6c5449546a3e51fe36981779d8b/e2e/f50661ca-a4ea-4f47-bb91-8d6c76b170cf.md"
$docName1 = "e25f3909-71f0-446d-a71c-4e2ddb72b79a.md"
$docName2 = "f50661ca-a4ea-4f47-bb91-8d6c76b170cf.md"

$hyperlinkColor = 15570276

# Positional-parameter function (named parameters are not reliably supported
# by this COM-interop runtime, so avoid "-param value" call syntax).
function Set-HandbackRow($ws, $row, $docUrl, $docName, $targetFileName, $handbackDateTime) {
    # Status: "In Translation" -> "Handed back: in sync with en-US"
    $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"

    # Latest Target File (column I) - hyperlink to the source doc, same as column A
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 9), $docUrl, "", "", $docName) | Out-Null
    $ws.Cells.Item($row, 9).Font.Underline = $true
    $ws.Cells.Item($row, 9).Font.Color = $hyperlinkColor

    # Latest Handback File (column J)
    $ws.Cells.Item($row, 10).Value = $targetFileName

    # Latest Handback DateTime (column K)
    $ws.Cells.Item($row, 11).Value = $handbackDateTime
}

# zh-cn sheet
Set-HandbackRow $ws2 2 $docUrl1 $docName1 "e25f3909-71f0-446d-a71c-4e2ddb72b79a.76efe022161d08dcd7722851807aa88ddc8b317e.zh-cn.xlf" "2016-08-18 02:20:29"
Set-HandbackRow $ws2 3 $docUrl2 $docName2 "f50661ca-a4ea-4f47-bb91-8d6c76b170cf.f3e904ac020f6fb00c84d5e7b1ce8c07d26db0b3.zh-cn.xlf" "2016-08-18 02:20:29"

# de-de sheet
Set-HandbackRow $ws3 2 $docUrl1 $docName1 "e25f3909-71f0-446d-a71c-4e2ddb72b79a.76efe022161d08dcd7722851807aa88ddc8b317e.de-de.xlf" "2016-08-18 02:20:36"
Set-HandbackRow $ws3 3 $docUrl2 $docName2 "f50661ca-a4ea-4f47-bb91-8d6c76b170cf.f3e904ac020f6fb00c84d5e7b1ce8c07d26db0b3.de-de.xlf" "2016-08-18 02:20:36"

# Overview sheet also rolls up the same "Status" string per language (columns E=zh-cn, F=de-de)
$ws1.Cells.Item(2, 5).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2, 6).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"

# Widen columns that now hold longer content (Status, Latest Target File, Latest Handback File)
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664
